$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell value updates (reflecting removed 'small_talk' row merged into
#     wellbeing_query/capabilities, row re-numbering, expanded phrase lists,
#     and the new AskAboutRestaurant intent row) ---
$ws.Range("B2").Value = 'cancel, exit, quit, bye, goodbye,close,shut down,no more,end,no more,log off'
$ws.Range("B4").Value = 'hello, hi, greetings, hey there,hey'
$ws.Range("B5").Value = 'how are you, how''s it going, how’s your day,what''s up, what’s happening, what’s new, whats up'
$ws.Range("B6").Value = 'what can you do, what can you help with, assist, help me,features'
$ws.Range("C6").Value = 'I can assist you with booking a table and answer questions., Type ''book'' to make a reservation or ask me anything and I’ll do my best to answer!,Not much just waiting to assist you!, Just here and ready to help!, I’m here feel free to ask me anything!'
$ws.Range("A7").Value = 'thanks'
$ws.Range("B7").Value = 'thank you, thanks, appreciate it,grateful,owe you,thankful'
$ws.Range("C7").Value = 'You''re very welcome! Anything else I can help with?, Glad to be of assistance!, It’s my pleasure, feel free to ask more questions!'
$ws.Range("A8").Value = 'positive_responses'
$ws.Range("B8").Value = 'y, yes, okay, confirm, ok, yeah, sure ,no problem,correct,grateful'
$ws.Range("C8").Value = 'Glad to be of assistance!'
$ws.Range("A9").Value = 'negative_responses'
$ws.Range("B9").Value = 'n, no, nah, nope,don''t agree,wrong,not sure,not interested,not up for it'
$ws.Range("C9").Value = 'I am sorry for that'
$ws.Range("A10").Value = 'name'
$ws.Range("B10").Value = 'my name, do you remember my name, do you know my name,  who am i, call my name, what''s my name, call me, my name is'
$ws.Range("A11").Value = 'change'
$ws.Range("B11").Value = 'change, change my name, update, update my name,a different name'
$ws.Range("A12").Value = 'time'
$ws.Range("B12").Value = 'time, current time, what time is it, do you know the current time,check the time,early or late.what hour'
$ws.Range("A13").Value = 'AskAboutRestaurant'
$ws.Range("B13").Value = 'restaurant, location, place, restaurant name, menu, food, dish, items, address, where, hours, operating hours, opening hours, opening time, offers, special offers, promotions, discounts,serving dinner,restaurant''s schedule'

# Row 10's Response cell is no longer used (AskAboutRestaurant has no Response)
$ws.Range("C10").Clear()

# --- Style fixes so per-cell formatting keeps tracking its intent row ---
$ws.Range("A8:B8").ClearFormats()
$ws.Range("A2").Copy()
$ws.Range("B9").PasteSpecial(-4122)
$ws.Range("A10").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Sheet view / column width adjustments ---
$ws.Columns("A").ColumnWidth = 51.25
$ws.Columns("B").ColumnWidth = 170.375
$ws.Columns("C").ColumnWidth = 206.125

$excel.ActiveWindow.Zoom = 130
